$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows for "Gruppe 1" (row 2) and "Gruppe 13" (row 14).
# Delete bottom-up so row indices above the deleted row are unaffected.
$ws.Rows(14).Delete()
$ws.Rows(2).Delete()

# Update remaining rows with new song-titled group names and updated image paths.
$ws.Range("A2").Value = "Gruppe 2 - UMBRELLA"
$ws.Range("B2").Value = "/album_covers/gruppe_2.png"

$ws.Range("A3").Value = "Gruppe 3 - I LOVE IT"
$ws.Range("B3").Value = "/album_covers/gruppe_3.png"

$ws.Range("A4").Value = "Gruppe 4 - SLIM SHADY"
$ws.Range("B4").Value = "/album_covers/gruppe_4.png"

$ws.Range("A5").Value = "Gruppe 5 - THAT’S WHAT IT IS"
$ws.Range("B5").Value = "/album_covers/artist_5.jpg"

$ws.Range("A6").Value = "Gruppe 6 - PHOENIX"
$ws.Range("B6").Value = "/album_covers/gruppe_6.png"

$ws.Range("A7").Value = "Gruppe 7 - AMERICAN IDIOT"
$ws.Range("B7").Value = "/album_covers/artist_7.jpg"

$ws.Range("A8").Value = "Gruppe 8 - LOVE YOU LIKE A LOVE SONG"
$ws.Range("B8").Value = "/album_covers/gruppe_8.png"

$ws.Range("A9").Value = "Gruppe 9 - THE DINER"
$ws.Range("B9").Value = "/album_covers/gruppe_9.png"

$ws.Range("A10").Value = "Gruppe 10 - MOVES LIKE JAGGER"
$ws.Range("B10").Value = "/album_covers/gruppe_10.png"

$ws.Range("A11").Value = "Gruppe 11 - NOW OR NEVER"
$ws.Range("B11").Value = "/album_covers/gruppe_11.png"

$ws.Range("A12").Value = "Gruppe 12 - ON THE FLOOR"
$ws.Range("B12").Value = "/album_covers/gruppe_12.png"

$ws.Range("A13").Value = "Gruppe 14 - UPTOWN FUNK"
$ws.Range("B13").Value = "/album_covers/artist_14.jpg"

$ws.Range("A14").Value = "Gruppe 15 - ROCKY"
$ws.Range("B14").Value = "/album_covers/gruppe_15.png"

$ws.Range("A15").Value = "Gruppe 16 - I’M STILL STANDING"
$ws.Range("B15").Value = "/album_covers/gruppe_16.jpg"

$ws.Range("A16").Value = "Gruppe 17 - HOTEL ROOM SERVICE"
$ws.Range("B16").Value = "/album_covers/artist_17.jpg"

$ws.Range("A17").Value = "Gruppe 18 - HER"
$ws.Range("B17").Value = "/album_covers/gruppe_18.jpg"

$ws.Range("A18").Value = "Gruppe 19 - DRACULA"
$ws.Range("B18").Value = "/album_covers/artist_19.jpg"

# Widen column A to fit the new, longer group/song names.
$ws.Columns(1).ColumnWidth = 60.69921875

# Move selection to the new last data row and park scroll at the top.
$ws.Range("B19").Select()
